$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "44.083.31", "239.79") are preserved exactly as text, matching
# the original inline-string cell contents rather than being coerced
# into floating point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.083.31"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.357.31"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.683"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.79"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.68%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.94"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.03%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.11"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.36"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +15.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.31"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.09%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.714.93"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.54"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.910"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.362.24"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.939.05"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000102"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.48%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.18"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.25"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.40%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("B25").Value = "ImmutableX"

$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.86"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +17.73%  "

$ws.Range("B26").Value = "WEMIXToken"

$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.72"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.35%  "

$ws.Range("B27").Value = "PancakeSwap"

$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.51"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.76"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.26%  "

$ws.Range("B29").Value = "Toncoin"

$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.61%  "

$ws.Range("B30").Value = "EthereumClassic"

$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.90"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.91"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.37%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.137"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0762"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.31"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.38"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.73"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.37"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.36"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0283"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.210"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +12.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.48"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.109"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +11.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.14"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.67%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.78"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +7.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.53"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +9.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.25"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.51%  "

$ws.Range("B49").Value = "Aave"

$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.39"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.36%  "

$ws.Range("B50").Value = "ARBITRUM"

$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.18"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.66"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +9.15%  "
